$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-05-22 Monday", $false, $false, $false, $false, $false, $true, 1, $false, "2023-05-23 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("2+71=73", $false, $false, $false, $false, $false, $true, 1, $false, "77-37=40", 2) | Out-Null
$d.Content.Find.Execute("53-46=7", $false, $false, $false, $false, $false, $true, 1, $false, "71-51=20", 2) | Out-Null
$d.Content.Find.Execute("51-21=30", $false, $false, $false, $false, $false, $true, 1, $false, "92-40=52", 2) | Out-Null
$d.Content.Find.Execute("21+35=56", $false, $false, $false, $false, $false, $true, 1, $false, "15+46=61", 2) | Out-Null
$d.Content.Find.Execute("87+2=89", $false, $false, $false, $false, $false, $true, 1, $false, "15+65=80", 2) | Out-Null
$d.Content.Find.Execute("71+18=89", $false, $false, $false, $false, $false, $true, 1, $false, "70-7=63", 2) | Out-Null
$d.Content.Find.Execute("15-13=2", $false, $false, $false, $false, $false, $true, 1, $false, "44+20=64", 2) | Out-Null
$d.Content.Find.Execute("80-47=33", $false, $false, $false, $false, $false, $true, 1, $false, "58+18=76", 2) | Out-Null
$d.Content.Find.Execute("80-36=44", $false, $false, $false, $false, $false, $true, 1, $false, "11+74=85", 2) | Out-Null
$d.Content.Find.Execute("18-5=13", $false, $false, $false, $false, $false, $true, 1, $false, "82-22=60", 2) | Out-Null
$d.Content.Find.Execute("23-16=7", $false, $false, $false, $false, $false, $true, 1, $false, "42-21=21", 2) | Out-Null
$d.Content.Find.Execute("51-31=20", $false, $false, $false, $false, $false, $true, 1, $false, "23+62=85", 2) | Out-Null
$d.Content.Find.Execute("82-60=22", $false, $false, $false, $false, $false, $true, 1, $false, "55-24=31", 2) | Out-Null
$d.Content.Find.Execute("46+43=89", $false, $false, $false, $false, $false, $true, 1, $false, "41-30=11", 2) | Out-Null
$d.Content.Find.Execute("15-7=8", $false, $false, $false, $false, $false, $true, 1, $false, "66+13=79", 2) | Out-Null
$d.Content.Find.Execute("44+43=87", $false, $false, $false, $false, $false, $true, 1, $false, "74-16=58", 2) | Out-Null
$d.Content.Find.Execute("25+72=97", $false, $false, $false, $false, $false, $true, 1, $false, "37-21=16", 2) | Out-Null
$d.Content.Find.Execute("25+42=67", $false, $false, $false, $false, $false, $true, 1, $false, "85-55=30", 2) | Out-Null
$d.Content.Find.Execute("4+26=30", $false, $false, $false, $false, $false, $true, 1, $false, "13+61=74", 2) | Out-Null
$d.Content.Find.Execute("11+23=34", $false, $false, $false, $false, $false, $true, 1, $false, "36-0=36", 2) | Out-Null
$d.Content.Find.Execute("74-28=46", $false, $false, $false, $false, $false, $true, 1, $false, "41+11=52", 2) | Out-Null
$d.Content.Find.Execute("79-63=16", $false, $false, $false, $false, $false, $true, 1, $false, "91-23=68", 2) | Out-Null
$d.Content.Find.Execute("24+37=61", $false, $false, $false, $false, $false, $true, 1, $false, "54-8=46", 2) | Out-Null
$d.Content.Find.Execute("38+10=48", $false, $false, $false, $false, $false, $true, 1, $false, "25+39=64", 2) | Out-Null
$d.Content.Find.Execute("86-67=19", $false, $false, $false, $false, $false, $true, 1, $false, "76-11=65", 2) | Out-Null
$d.Content.Find.Execute("73-54=19", $false, $false, $false, $false, $false, $true, 1, $false, "12+12=24", 2) | Out-Null
$d.Content.Find.Execute("20-16=4", $false, $false, $false, $false, $false, $true, 1, $false, "69+1=70", 2) | Out-Null
$d.Content.Find.Execute("94-91=3", $false, $false, $false, $false, $false, $true, 1, $false, "38+4=42", 2) | Out-Null
$d.Content.Find.Execute("36+13=49", $false, $false, $false, $false, $false, $true, 1, $false, "12+67=79", 2) | Out-Null
$d.Content.Find.Execute("3-1=2", $false, $false, $false, $false, $false, $true, 1, $false, "53+35=88", 2) | Out-Null
$d.Content.Find.Execute("36+28=64", $false, $false, $false, $false, $false, $true, 1, $false, "44-14=30", 2) | Out-Null
$d.Content.Find.Execute("80-17=63", $false, $false, $false, $false, $false, $true, 1, $false, "35+16=51", 2) | Out-Null
$d.Content.Find.Execute("71-22=49", $false, $false, $false, $false, $false, $true, 1, $false, "83-33=50", 2) | Out-Null
$d.Content.Find.Execute("87-40=47", $false, $false, $false, $false, $false, $true, 1, $false, "70-47=23", 2) | Out-Null
$d.Content.Find.Execute("18+80=98", $false, $false, $false, $false, $false, $true, 1, $false, "22+57=79", 2) | Out-Null
$d.Content.Find.Execute("10+43=53", $false, $false, $false, $false, $false, $true, 1, $false, "97-4=93", 2) | Out-Null
$d.Content.Find.Execute("14+1=15", $false, $false, $false, $false, $false, $true, 1, $false, "61-5=56", 2) | Out-Null
$d.Content.Find.Execute("34+43=77", $false, $false, $false, $false, $false, $true, 1, $false, "28-20=8", 2) | Out-Null
$d.Content.Find.Execute("25+51=76", $false, $false, $false, $false, $false, $true, 1, $false, "69-0=69", 2) | Out-Null
$d.Content.Find.Execute("37+41=78", $false, $false, $false, $false, $false, $true, 1, $false, "32+51=83", 2) | Out-Null
$d.Content.Find.Execute("63-6=57", $false, $false, $false, $false, $false, $true, 1, $false, "45-25=20", 2) | Out-Null
$d.Content.Find.Execute("74+1=75", $false, $false, $false, $false, $false, $true, 1, $false, "15+41=56", 2) | Out-Null
$d.Content.Find.Execute("87-86=1", $false, $false, $false, $false, $false, $true, 1, $false, "36+48=84", 2) | Out-Null
$d.Content.Find.Execute("68+17=85", $false, $false, $false, $false, $false, $true, 1, $false, "99-88=11", 2) | Out-Null
$d.Content.Find.Execute("55-50=5", $false, $false, $false, $false, $false, $true, 1, $false, "1+13=14", 2) | Out-Null
$d.Content.Find.Execute("95-11=84", $false, $false, $false, $false, $false, $true, 1, $false, "0+0=0", 2) | Out-Null
$d.Content.Find.Execute("12+25=37", $false, $false, $false, $false, $false, $true, 1, $false, "62-46=16", 2) | Out-Null
$d.Content.Find.Execute("42+7=49", $false, $false, $false, $false, $false, $true, 1, $false, "6+70=76", 2) | Out-Null
$d.Content.Find.Execute("70-43=27", $false, $false, $false, $false, $false, $true, 1, $false, "71-31=40", 2) | Out-Null
$d.Content.Find.Execute("47-23=24", $false, $false, $false, $false, $false, $true, 1, $false, "67-17=50", 2) | Out-Null
$d.Content.Find.Execute("68-27=41", $false, $false, $false, $false, $false, $true, 1, $false, "0+41=41", 2) | Out-Null
$d.Content.Find.Execute("55-21=34", $false, $false, $false, $false, $false, $true, 1, $false, "41+9=50", 2) | Out-Null
$d.Content.Find.Execute("28+20=48", $false, $false, $false, $false, $false, $true, 1, $false, "89-45=44", 2) | Out-Null
$d.Content.Find.Execute("20+28=48", $false, $false, $false, $false, $false, $true, 1, $false, "98-85=13", 2) | Out-Null
$d.Content.Find.Execute("64-58=6", $false, $false, $false, $false, $false, $true, 1, $false, "26-8=18", 2) | Out-Null
$d.Content.Find.Execute("67-7=60", $false, $false, $false, $false, $false, $true, 1, $false, "34+36=70", 2) | Out-Null
$d.Content.Find.Execute("42+13=55", $false, $false, $false, $false, $false, $true, 1, $false, "94-3=91", 2) | Out-Null
$d.Content.Find.Execute("99-74=25", $false, $false, $false, $false, $false, $true, 1, $false, "20+33=53", 2) | Out-Null
$d.Content.Find.Execute("28+1=29", $false, $false, $false, $false, $false, $true, 1, $false, "34+3=37", 2) | Out-Null
$d.Content.Find.Execute("79-53=26", $false, $false, $false, $false, $false, $true, 1, $false, "28-5=23", 2) | Out-Null
$d.Content.Find.Execute("27+38=65", $false, $false, $false, $false, $false, $true, 1, $false, "13+70=83", 2) | Out-Null
$d.Content.Find.Execute("94-21=73", $false, $false, $false, $false, $false, $true, 1, $false, "16+54=70", 2) | Out-Null
$d.Content.Find.Execute("40-36=4", $false, $false, $false, $false, $false, $true, 1, $false, "41-13=28", 2) | Out-Null
$d.Content.Find.Execute("75-25=50", $false, $false, $false, $false, $false, $true, 1, $false, "27+48=75", 2) | Out-Null
$d.Content.Find.Execute("83-79=4", $false, $false, $false, $false, $false, $true, 1, $false, "95-76=19", 2) | Out-Null
$d.Content.Find.Execute("87-78=9", $false, $false, $false, $false, $false, $true, 1, $false, "32+13=45", 2) | Out-Null
$d.Content.Find.Execute("51-23=28", $false, $false, $false, $false, $false, $true, 1, $false, "99-86=13", 2) | Out-Null
$d.Content.Find.Execute("39+17=56", $false, $false, $false, $false, $false, $true, 1, $false, "10-3=7", 2) | Out-Null
$d.Content.Find.Execute("77-32=45", $false, $false, $false, $false, $false, $true, 1, $false, "59-22=37", 2) | Out-Null
$d.Content.Find.Execute("68-31=37", $false, $false, $false, $false, $false, $true, 1, $false, "31+38=69", 2) | Out-Null
$d.Content.Find.Execute("56-20=36", $false, $false, $false, $false, $false, $true, 1, $false, "92-80=12", 2) | Out-Null
$d.Content.Find.Execute("7+41=48", $false, $false, $false, $false, $false, $true, 1, $false, "37+19=56", 2) | Out-Null
$d.Content.Find.Execute("12+1=13", $false, $false, $false, $false, $false, $true, 1, $false, "5+93=98", 2) | Out-Null
$d.Content.Find.Execute("99-68=31", $false, $false, $false, $false, $false, $true, 1, $false, "66-38=28", 2) | Out-Null
$d.Content.Find.Execute("40+26=66", $false, $false, $false, $false, $false, $true, 1, $false, "74+13=87", 2) | Out-Null
$d.Content.Find.Execute("76-33=43", $false, $false, $false, $false, $false, $true, 1, $false, "99-52=47", 2) | Out-Null
$d.Content.Find.Execute("41+36=77", $false, $false, $false, $false, $false, $true, 1, $false, "77-15=62", 2) | Out-Null
$d.Content.Find.Execute("71+26=97", $false, $false, $false, $false, $false, $true, 1, $false, "55-28=27", 2) | Out-Null
$d.Content.Find.Execute("71+20=91", $false, $false, $false, $false, $false, $true, 1, $false, "44-9=35", 2) | Out-Null
$d.Content.Find.Execute("96-37=59", $false, $false, $false, $false, $false, $true, 1, $false, "95-4=91", 2) | Out-Null
$d.Content.Find.Execute("34+15=49", $false, $false, $false, $false, $false, $true, 1, $false, "89-6=83", 2) | Out-Null
$d.Content.Find.Execute("59-52=7", $false, $false, $false, $false, $false, $true, 1, $false, "30+62=92", 2) | Out-Null
$d.Content.Find.Execute("45-37=8", $false, $false, $false, $false, $false, $true, 1, $false, "7+89=96", 2) | Out-Null
$d.Content.Find.Execute("89+2=91", $false, $false, $false, $false, $false, $true, 1, $false, "44-24=20", 2) | Out-Null
$d.Content.Find.Execute("76-14=62", $false, $false, $false, $false, $false, $true, 1, $false, "57-1=56", 2) | Out-Null
$d.Content.Find.Execute("66+24=90", $false, $false, $false, $false, $false, $true, 1, $false, "11-0=11", 2) | Out-Null
$d.Content.Find.Execute("96-84=12", $false, $false, $false, $false, $false, $true, 1, $false, "82-63=19", 2) | Out-Null
$d.Content.Find.Execute("93-56=37", $false, $false, $false, $false, $false, $true, 1, $false, "24+21=45", 2) | Out-Null
$d.Content.Find.Execute("20-9=11", $false, $false, $false, $false, $false, $true, 1, $false, "57-7=50", 2) | Out-Null
$d.Content.Find.Execute("63-9=54", $false, $false, $false, $false, $false, $true, 1, $false, "13+9=22", 2) | Out-Null
$d.Content.Find.Execute("59-25=34", $false, $false, $false, $false, $false, $true, 1, $false, "88-33=55", 2) | Out-Null
$d.Content.Find.Execute("45-29=16", $false, $false, $false, $false, $false, $true, 1, $false, "21-20=1", 2) | Out-Null
$d.Content.Find.Execute("6+83=89", $false, $false, $false, $false, $false, $true, 1, $false, "75-40=35", 2) | Out-Null
$d.Content.Find.Execute("41-16=25", $false, $false, $false, $false, $false, $true, 1, $false, "35-17=18", 2) | Out-Null
$d.Content.Find.Execute("90-36=54", $false, $false, $false, $false, $false, $true, 1, $false, "23-10=13", 2) | Out-Null
$d.Content.Find.Execute("4+85=89", $false, $false, $false, $false, $false, $true, 1, $false, "59+34=93", 2) | Out-Null
$d.Content.Find.Execute("83-15=68", $false, $false, $false, $false, $false, $true, 1, $false, "40+57=97", 2) | Out-Null
$d.Content.Find.Execute("48+28=76", $false, $false, $false, $false, $false, $true, 1, $false, "87-6=81", 2) | Out-Null
$d.Content.Find.Execute("13+79=92", $false, $false, $false, $false, $false, $true, 1, $false, "82-35=47", 2) | Out-Null
$d.Content.Find.Execute("53-14=39", $false, $false, $false, $false, $false, $true, 1, $false, "95-21=74", 2) | Out-Null
